$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team 2")
$ws.Activate()

# --- Team number cell: centre align ---
$ws.Range("B1").HorizontalAlignment = -4108

# --- Make room: insert 3 blank rows after row 11 (old rows 12-17 shift to 15-20) ---
$ws.Range("A12:A14").EntireRow.Insert()

# --- Row 11: Start (Amsterdam time) ---
$ws.Range("A11").Value = "Start (Amsterdam time):"
$ws.Range("B11").Value = "09:00 Saturday"

# --- Row 12: End (Amsterdam time) [NEW] ---
$ws.Range("A12").Value = "End (Amsterdam time):"
$ws.Range("B12").Value = "09:30 Sunday (Deadline was 10:00 due to Summer time shift)"

# --- Row 13: Results: ---
$ws.Range("A13").Value = "Results:"
$ws.Range("B13").Value = ""
$ws.Range("B13").ClearFormats()

# --- Row 14: Everything is complete under: [NEW] ---
$ws.Range("A14").Value = "Everything is complete under:"

# --- Row 15: Code: ---
$ws.Range("A15").Value = "Code:"

# --- Row 16: README: ---
$ws.Range("A16").Value = "README:"

# --- Row 17: One pager: ---
$ws.Range("A17").Value = "One pager:"

# --- Row 18: Pitch ppt: ---
$ws.Range("A18").Value = "Pitch ppt:"

# --- Row 19: Pitch video: ---
$ws.Range("A19").Value = "Pitch video:"

# --- Row 20: CODE DEMO video: [NEW] ---
$ws.Range("A20").Value = "CODE DEMO video:"

# ============ Styling ============

# A11 & A12: wrap text + header fill
$ws.Range("A11:A12").WrapText = $true
$ws.Range("A11:A12").Interior.ThemeColor = 4

# B11 & B12: time number format + header fill
$ws.Range("B11:B12").NumberFormat = "h:mm"
$ws.Range("B11:B12").Interior.ThemeColor = 4

# A14: bold (same as "Results:"/"Code:" headers) + header fill
$ws.Range("A14").Font.Bold = $true
$ws.Range("A14").Interior.ThemeColor = 4

# B14: hyperlink (display text == address, engine still records it) + header fill
$ws.Hyperlinks.Add($ws.Range("B14"), "https://github.com/bsirmacek/SchoolofAI_Healthcare_Hackathon_the_Netherlands/tree/master/Team2_Lotad.AI_completed") | Out-Null
$ws.Range("B14").Interior.ThemeColor = 4

# A15:A20: header fill
$ws.Range("A15:A20").Interior.ThemeColor = 4

# B15:B18: hyperlink with explicit display text (trailing space kept) + header fill + vertical centre
$ws.Hyperlinks.Add($ws.Range("B15"), "https://surfdrive.surf.nl/files/index.php/s/XDqe1jx5GETaf4X", "", "", "https://surfdrive.surf.nl/files/index.php/s/XDqe1jx5GETaf4X ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B16"), "https://surfdrive.surf.nl/files/index.php/s/XDqe1jx5GETaf4X", "", "", "https://surfdrive.surf.nl/files/index.php/s/XDqe1jx5GETaf4X ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B17"), "https://surfdrive.surf.nl/files/index.php/s/XDqe1jx5GETaf4X", "", "", "https://surfdrive.surf.nl/files/index.php/s/XDqe1jx5GETaf4X ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B18"), "https://surfdrive.surf.nl/files/index.php/s/XDqe1jx5GETaf4X", "", "", "https://surfdrive.surf.nl/files/index.php/s/XDqe1jx5GETaf4X ") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B19"), "https://youtu.be/wijWW5HopG0") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B20"), "https://youtu.be/ZnxV_MA78os") | Out-Null

$ws.Range("B15:B20").Interior.ThemeColor = 4
$ws.Range("B15:B20").VerticalAlignment = -4108

# ============ Sheet-level bits ============
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A1:A5").Select()
